$wb = $excel.ActiveWorkbook

# Update the "Test Format" for the C / Assignment 3 row on Sheet1:
# it was "Unit", now it is "Output" (matches sharedStrings entry already
# used elsewhere in the workbook, e.g. C5).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C4").Value = "Output"

# Switch the active/selected sheet from "A4" back to "Sheet1", and update
# Sheet1's selected cell to F25 (previously K32). Activating the sheet
# also clears the old tabSelected flag on whichever sheet was active
# before (the "A4" sheet), matching its selection staying at F13.
$ws1.Activate()
$ws1.Range("F25").Select() | Out-Null
